# "Firestore write ok, reading still missing"
#
# Extends the Tabelle1 documentation table:
#   - row 23 (areaMessages) now also gets a "me_button_send" action element
#   - three new rows (24-26) document further "Messages" area elements
#   - rows 27-31 document a brand-new "Navigation Items" area
#
# New rows re-use the same look (font) as the rest of the "Action Elements"
# column, i.e. the style already used by cells such as B6/B21 (cellXfs
# index 1: 10pt "Open Sans", color #333333). We copy that formatting from
# an existing styled cell via Copy/PasteSpecial(xlPasteFormats) instead of
# touching .Font.* directly, since the latter creates redundant
# font/style entries instead of reusing the existing one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- formatting -----------------------------------------------------
# B23, C23..C25, B24..B26, C27..C31 all carry the existing "styled" look
# (s="1"); B27 ("Navigation Items") stays unstyled, just like B15/B20 etc.
$ws.Range("B6").Copy()
$styledTargets = @("B23", "C23", "B24", "C24", "B25", "C25", "B26", "C27", "C28", "C29", "C30", "C31")
foreach ($target in $styledTargets) {
  $ws.Range($target).PasteSpecial(-4122)
}

# --- values -----------------------------------------------------------
# Written in the same order as the source edit so the shared-strings table
# comes out in the same sequence.
$ws.Range("B27").Value = "Navigation Items"
$ws.Range("C29").Value = "nav_sign_out"
$ws.Range("C30").Value = "nav_sign_in"
$ws.Range("C31").Value = "nav_account"
$ws.Range("C27").Value = "nav_pluto23"
$ws.Range("C28").Value = "nav_home"
$ws.Range("C24").Value = "me_div_messages"
$ws.Range("C25").Value = "me_body_text"
$ws.Range("C23").Value = "me_button_send"
$ws.Range("B23").Value = "areaMessages"

# Matches the saved selection in the edited workbook.
$ws.Range("C23").Select()
